$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3060
$ws1.Range("F4").Value = 56
$ws1.Range("F5").Value = 41
$ws1.Range("F7").Value = 1048
$ws1.Range("F8").Value = 14783
$ws1.Range("F10").Value = 136
$ws1.Range("F11").Value = 5894
$ws1.Range("F12").Value = 603
$ws1.Range("F13").Value = 84
$ws1.Range("F14").Value = 49
$ws1.Range("F15").Value = 81
$ws1.Range("F16").Value = 1243
$ws1.Range("F18").Value = 93
$ws1.Range("F19").Value = 196
$ws1.Range("F20").Value = 810
$ws1.Range("F21").Value = 2947
$ws1.Range("F22").Value = 96
$ws1.Range("F23").Value = 10698
$ws1.Range("F24").Value = 1207
$ws1.Range("F26").Value = 110
$ws1.Range("F28").Value = 250

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3060
$ws4.Range("F5").Value = 56
$ws4.Range("F6").Value = 41
$ws4.Range("F8").Value = 1048
$ws4.Range("F9").Value = 14783
$ws4.Range("F11").Value = 136
$ws4.Range("F12").Value = 5894
$ws4.Range("F13").Value = 603
$ws4.Range("F14").Value = 84
$ws4.Range("F15").Value = 49
$ws4.Range("F16").Value = 81
$ws4.Range("F17").Value = 1243
$ws4.Range("F19").Value = 93
$ws4.Range("F20").Value = 196
$ws4.Range("F21").Value = 810
$ws4.Range("F22").Value = 2947
$ws4.Range("F23").Value = 96
$ws4.Range("F25").Value = 10698
$ws4.Range("F26").Value = 1207
$ws4.Range("F28").Value = 110
$ws4.Range("F30").Value = 250

$wb.Save()
